# Apply November 2018 Unofficial Model Run updates (CBRFC forecast + model output refresh)
# to the MtomToCrss_Annual Trace sheets. Values below are the refreshed model outputs.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("Trace3")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8854562.1209468655
$ws.Range("H3").Value = 8998832.8250046521

$ws = $wb.Worksheets("Trace4")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8477489.750057783

$ws = $wb.Worksheets("Trace5")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("H3").Value = 9322938.8360353373

$ws = $wb.Worksheets("Trace6")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8544080.4700577836

$ws = $wb.Worksheets("Trace7")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("H3").Value = 9253869.5300353374

$ws = $wb.Worksheets("Trace8")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("H3").Value = 9182721.3480353393

$ws = $wb.Worksheets("Trace9")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8999999.9999775533
$ws.Range("H3").Value = 9220479.9720353391

$ws = $wb.Worksheets("Trace10")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8452673.5360577814

$ws = $wb.Worksheets("Trace11")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8564528.281613294
$ws.Range("H3").Value = 8745731.51567108

$ws = $wb.Worksheets("Trace12")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("C3").Value = 263899.52890398417
$ws.Range("E3").Value = 377046.52890402172
$ws.Range("F3").Value = 8366132.5451393249
$ws.Range("H3").Value = 8539234.8991971128

$ws = $wb.Worksheets("Trace13")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8449825.4580577817

$ws = $wb.Worksheets("Trace14")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8510853.7240577824

$ws = $wb.Worksheets("Trace15")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8671869.3440577816

$ws = $wb.Worksheets("Trace16")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8999999.9999775533
$ws.Range("H3").Value = 9171027.6120353378

$ws = $wb.Worksheets("Trace17")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8550902.8760577813

$ws = $wb.Worksheets("Trace18")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("H3").Value = 9166368.9920353387

$ws = $wb.Worksheets("Trace19")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("H3").Value = 9105269.704035338

$ws = $wb.Worksheets("Trace20")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8999999.9999775533
$ws.Range("H3").Value = 8968184.204035338

$ws = $wb.Worksheets("Trace21")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("C3").Value = 919439.52890378458
$ws.Range("E3").Value = 1032586.5289038221
$ws.Range("F3").Value = 8626535.657290794
$ws.Range("H3").Value = 8855643.4713485781

$ws = $wb.Worksheets("Trace22")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8999999.9999775533
$ws.Range("H3").Value = 9179659.4840353373

$ws = $wb.Worksheets("Trace23")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8373740.0240577832

$ws = $wb.Worksheets("Trace24")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8384084.7640577815

$ws = $wb.Worksheets("Trace25")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8404937.7940577827

$ws = $wb.Worksheets("Trace26")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8772637.3478710875
$ws.Range("H3").Value = 8917508.7719288729

$ws = $wb.Worksheets("Trace27")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8409652.7440577839

$ws = $wb.Worksheets("Trace28")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8999999.9999775533
$ws.Range("H3").Value = 9295759.3240353372

$ws = $wb.Worksheets("Trace29")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8999999.9999775533
$ws.Range("H3").Value = 9196913.5740353391

$ws = $wb.Worksheets("Trace30")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8455034.3440577816

$ws = $wb.Worksheets("Trace31")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("H3").Value = 9196925.9240353368

$ws = $wb.Worksheets("Trace32")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8999999.9999775533
$ws.Range("H3").Value = 9219887.1840353366

$ws = $wb.Worksheets("Trace33")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8511829.9610227067
$ws.Range("H3").Value = 8767557.2750804927

$ws = $wb.Worksheets("Trace34")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8999999.9999775533
$ws.Range("H3").Value = 9087757.4640353378

$ws = $wb.Worksheets("Trace35")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8314596.0640577823

$ws = $wb.Worksheets("Trace36")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8365743.0240577832

$ws = $wb.Worksheets("Trace37")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("F3").Value = 8737151.7454127055
$ws.Range("H3").Value = 8908334.2294704895

$ws = $wb.Worksheets("Trace38")
$ws.Range("H2").Value = 9315862.6005992647
$ws.Range("C3").Value = 75719.528904368795
$ws.Range("E3").Value = 188866.52890440638
$ws.Range("F3").Value = 8229999.9999999981
$ws.Range("H3").Value = 8382461.3600577824
